# Update "想去人数" (number of people interested) counts for three events
# on the "展览" sheet and the duplicated rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 7649
$wsExhibition.Range("F6").Value = 4314
$wsExhibition.Range("F11").Value = 161

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 7649
$wsAllTypes.Range("F7").Value = 4314
$wsAllTypes.Range("F13").Value = 161
